$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2").Delete()
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "4873658"
$ws.Range("V1").Value = "Resultado"
$ws.Range("V2").Value = "Se han encontrado errores en la Validacion de la Propuesta"
$ws.Range("V2").Select()
